# Apply cryptos list updates to match target revision
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '68.506.06'
$ws.Cells.Item(2, 5).Value = '  +1.25%  '
$ws.Cells.Item(3, 4).Value = '2.648.88'
$ws.Cells.Item(3, 5).Value = '  +1.30%  '
$ws.Cells.Item(4, 5).Value = '  -0.05%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '599.50'
$ws.Cells.Item(5, 5).Value = '  +0.85%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '154.70'
$ws.Cells.Item(6, 5).Value = '  +1.66%  '
$ws.Cells.Item(7, 5).Value = '  -0.07%  '
$ws.Cells.Item(8, 5).Value = '  +0.70%  '
$ws.Cells.Item(9, 4).Value = '2.647.82'
$ws.Cells.Item(9, 5).Value = '  +1.31%  '
$ws.Cells.Item(10, 5).Value = '  +9.61%  '
$ws.Cells.Item(11, 5).Value = '  -0.19%  '
$ws.Cells.Item(12, 5).Value = '  +1.57%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.355'
$ws.Cells.Item(13, 5).Value = '  +2.51%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '28.32'
$ws.Cells.Item(14, 5).Value = '  +3.10%  '
$ws.Cells.Item(15, 5).Value = '  +3.37%  '
$ws.Cells.Item(16, 4).Value = '3.130.53'
$ws.Cells.Item(17, 4).Value = '68.455.94'
$ws.Cells.Item(17, 5).Value = '  +1.17%  '
$ws.Cells.Item(18, 4).Value = '2.641.99'
$ws.Cells.Item(18, 5).Value = '  +1.22%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '11.50'
$ws.Cells.Item(19, 5).Value = '  +2.83%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '366.74'
$ws.Cells.Item(20, 5).Value = '  -1.28%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '7.52'
$ws.Cells.Item(21, 5).Value = '  +1.47%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '4.41'
$ws.Cells.Item(22, 5).Value = '  +4.55%  '
$ws.Cells.Item(23, 5).Value = '  +2.46%  '
$ws.Cells.Item(24, 5).Value = '  +2.56%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '73.79'
$ws.Cells.Item(25, 5).Value = '  +1.66%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '0.999'
$ws.Cells.Item(26, 5).Value = '  -0.05%  '
$ws.Cells.Item(27, 5).Value = '  +0.99%  '
$ws.Cells.Item(28, 5).Value = '  +4.81%  '
$ws.Cells.Item(29, 4).Value = '2.780.52'
$ws.Cells.Item(29, 5).Value = '  +1.11%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '581.00'
$ws.Cells.Item(30, 5).Value = '  -2.31%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '0.999'
$ws.Cells.Item(31, 5).Value = '  -0.25%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '8.23'
$ws.Cells.Item(32, 5).Value = '  +5.68%  '
$ws.Cells.Item(33, 5).Value = '  +4.73%  '
$ws.Cells.Item(34, 5).Value = '  +1.72%  '
$ws.Cells.Item(35, 5).Value = '  +5.00%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '1.61'
$ws.Cells.Item(36, 5).Value = '  +6.49%  '
$ws.Cells.Item(37, 5).Value = '  -0.07%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '159.57'
$ws.Cells.Item(38, 5).Value = '  +0.95%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '19.51'
$ws.Cells.Item(39, 5).Value = '  +2.13%  '
$ws.Cells.Item(40, 2).Value = 'Stacks'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '1.91'
$ws.Cells.Item(40, 5).Value = '  +1.42%  '
$ws.Cells.Item(41, 2).Value = 'PolygonEcosystemToken'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.376'
$ws.Cells.Item(41, 5).Value = '  +2.39%  '
$ws.Cells.Item(42, 5).Value = '  +4.06%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '2.70'
$ws.Cells.Item(43, 5).Value = '  +1.04%  '
$ws.Cells.Item(44, 4).Value = '0.0₆0336'
$ws.Cells.Item(44, 5).Value = '  +13.94%  '
$ws.Cells.Item(45, 5).Value = '  +3.55%  '
$ws.Cells.Item(46, 5).Value = '  -0.01%  '
$ws.Cells.Item(47, 5).Value = '  +0.39%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '157.90'
$ws.Cells.Item(48, 5).Value = '  +1.24%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '3.80'
$ws.Cells.Item(49, 5).Value = '  +3.67%  '
$ws.Cells.Item(50, 5).Value = '  +2.63%  '
$ws.Cells.Item(51, 5).Value = '  +4.03%  '
